# Deploying to gh-pages from @ 03fb76b9ac146585df05395f980ae353fb99e762
# Updates the "8.7.1 child labour" sheet: fixes the sex/urbanisation
# header & value wording (plural/title-case forms) across the three
# language columns, capitalises several English category labels, and
# fills in the previously-empty "Functional difficulties in a child"
# row label in the Kyrgyz column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - "by sex" section header
$ws.Range("A6").Value = "Жынысы боюнча"
$ws.Range("B6").Value = "По полу"
$ws.Range("C6").Value = "By sex"

# Row 7 - male
$ws.Range("A7").Value = "Эркектер"
$ws.Range("B7").Value = "Мужчины"
$ws.Range("C7").Value = "Men"

# Row 8 - female
$ws.Range("A8").Value = "Аялдар"
$ws.Range("B8").Value = "Женщины"
$ws.Range("C8").Value = "Woman"

# Row 10 - urban
$ws.Range("A10").Value = "Шаар"
$ws.Range("C10").Value = "Urban"

# Row 11 - rural
$ws.Range("A11").Value = "Айыл"
$ws.Range("C11").Value = "Rural"

# Row 28 - does not attend school (English column only)
$ws.Range("C28").Value = "Does not attend"

# Row 29 - mother's education header (English column only)
$ws.Range("C29").Value = "Educationof mother"

# Row 30 - preschool or not/primary (English column only)
$ws.Range("C30").Value = "Preschool or not /primary"

# Row 31 - basic general (English column only)
$ws.Range("C31").Value = "Basic general"

# Row 32 - average total (English column only)
$ws.Range("C32").Value = "Average total"

# Row 33 - vocational primary/secondary (English column only)
$ws.Range("C33").Value = "Vocational primary /secondary"

# Row 34 - higher (English column only)
$ws.Range("C34").Value = "Higher"

# Row 35 - "Functional difficulties in a child" header; the Kyrgyz cell
# (A35) was previously blank and now gets the matching label, picking up
# the same formatting already used by its row-mates B35/C35.
$ws.Range("A35").Font.Name = $ws.Range("B35").Font.Name
$ws.Range("A35").Font.Size = $ws.Range("B35").Font.Size
$ws.Range("A35").Font.Bold = $ws.Range("B35").Font.Bold
$ws.Range("A35").Font.Italic = $ws.Range("B35").Font.Italic
$ws.Range("A35").VerticalAlignment = $ws.Range("B35").VerticalAlignment
$ws.Range("A35").HorizontalAlignment = $ws.Range("B35").HorizontalAlignment
$ws.Range("A35").Value = "Баланын функционалдык кыйнчылыктары"

# Row 38 - wealth quintile header (English column only)
$ws.Range("C38").Value = "Wealth quintile"
